# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
#
# Column G (header "K") holds the strikeout count per outing. This re-derives
# the per-row K values (previously mis-populated) and writes the corrected
# values back into the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping of row -> corrected K value, per the regenerated save_data.
$kValues = @{
    2  = 2
    3  = 0
    4  = 1
    6  = 2
    7  = 0
    8  = 1
    9  = 2
    10 = 1
    11 = 0
    12 = 0
    13 = 2
    14 = 0
    15 = 2
    16 = 1
    17 = 0
    19 = 3
    20 = 1
    21 = 0
    22 = 2
    23 = 0
    24 = 2
    25 = 0
    26 = 2
    27 = 1
    28 = 1
    29 = 0
    30 = 1
    31 = 1
    32 = 3
    33 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
